# Fill in the four (still-blank) "answer" placeholders with the author's
# real answers, splitting the first one's content across two runs (matching
# how Word split the "Problem 4 took " / "about an hour." sentence), and
# drop the stray page-break-only paragraph that used to separate the last
# answer from the "Screenshots" section.

$d = $word.ActiveDocument

function Replace-AnswerInParagraph($paraIndex, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range.Duplicate()
    $rng.Find.Execute("answer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $rng.Text = $newText
}

# --- Q1: "How many hours did it take you to complete this assessment?" ---
# The answer paragraph is paragraph 5 ("A: answer"). The target splits the
# single "answer" run into two runs with identical run formatting, so we
# first clone the formatted run (Copy/Paste keeps its rPr, unlike
# InsertAfter which creates a run with no formatting at all), then give the
# clone a throwaway distinguishing format (Bold) so the two adjacent,
# equally-formatted runs aren't auto-merged while we edit their text, and
# finally clear that throwaway format again.
$p5 = $d.Paragraphs.Item(5)
$found = $p5.Range.Duplicate()
$found.Find.Execute("answer")
$found.Copy()
$found.Collapse(0)
$found.Paste()

$p5 = $d.Paragraphs.Item(5)
$run1 = $p5.Range.Duplicate()
$run1.Find.Execute("answer")
$run1.Bold = 1

$part1 = "Problem 1 took roughly two and a half to three hours. Problem 2 took an hour or so. Problem 3 took around four hours. Problem 4 took "
$part2 = "about an hour."

$run1.Text = $part1

$p5 = $d.Paragraphs.Item(5)
$run2 = $p5.Range.Duplicate()
$run2.Start = $run1.End
$run2.End = $p5.Range.End - 1
$run2.Text = $part2

$p5 = $d.Paragraphs.Item(5)
$run1fix = $p5.Range.Duplicate()
$run1fix.Start = $p5.Range.Start + 3
$run1fix.End = $run1fix.Start + $part1.Length
$run1fix.Bold = 0

# --- Q2: "What online resources did you use?" ---
Replace-AnswerInParagraph 9 "I used some Oracle documentation just to verify things about specific functions but that’s it. Maybe some Stack Overflow too."

# --- Q3: "Did you get help from any classmates?" ---
Replace-AnswerInParagraph 13 "No, I did not require any assistance from my peers for this endeavor."

# --- Q4: "Did you ask for help from an instructor?" ---
Replace-AnswerInParagraph 17 "Nay, I required not such a query."

# --- Q5: "Rate the difficulty of each problem..." ---
Replace-AnswerInParagraph 21 "Employing the use of the common “one-to-ten” scale, my dispositions are as follows: Problem 1 receives a 2/10. Problem 2 receives a 3/10. Problem 3 receives a 5/10. Problem 4 receives a 2/10. I am confident I can solve similar issues in the future."

# Remove the now-orphaned page-break-only paragraph that followed Q5's answer.
$pageBreakPara = $d.Paragraphs.Item(22)
$pageBreakPara.Range.Delete()

Write-Output "Applied answers and removed stray page break paragraph."
